$d = $word.ActiveDocument

$replacements = @(
    @{old = "79÷5="; new = "68÷2="},
    @{old = "50÷4="; new = "81÷8="},
    @{old = "54÷3="; new = "64÷3="},
    @{old = "40÷6="; new = "49÷2="},
    @{old = "54÷7="; new = "57÷7="},
    @{old = "47÷7="; new = "35÷8="},
    @{old = "32÷3="; new = "91÷4="},
    @{old = "38÷4="; new = "11÷9="},
    @{old = "27÷6="; new = "41÷9="},
    @{old = "75÷2="; new = "97÷8="},
    @{old = "59÷8="; new = "57÷4="},
    @{old = "26÷7="; new = "77÷5="},
    @{old = "47÷5="; new = "42÷6="},
    @{old = "66÷5="; new = "99÷3="},
    @{old = "74÷9="; new = "63÷9="},
    @{old = "91÷9="; new = "73÷9="},
    @{old = "46÷8="; new = "67÷4="},
    @{old = "89÷7="; new = "45÷6="},
    @{old = "28÷3="; new = "38÷3="},
    @{old = "44÷5="; new = "42÷9="},
    @{old = "41÷8="; new = "30÷3="},
    @{old = "38÷6="; new = "60÷4="},
    @{old = "58÷6="; new = "85÷9="},
    @{old = "51÷7="; new = "70÷8="},
    @{old = "78÷6="; new = "64÷9="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
